# "Partie bouton Achat et Terminer a continuer"
# Append two new client rows (30 and 31) to the bottom of the Clients table
# on the active sheet ("Feuil1"), mirroring the shape of the existing rows
# (22-29): card number in col A, client name in col B, bonus points in col C,
# credit balance in col D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30 : card "516122" for the existing client "Hana Murata" ---------
# Card numbers are stored as text (e.g. "111888", "222222", ...), even
# though they look numeric, so force a text format before typing the value
# in order to avoid Excel auto-converting it to a number. Restore the
# cell's style afterwards so no stray number-format is left behind.
$ws.Cells.Item(30, 1).NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "516122"
$ws.Cells.Item(30, 1).Style = "Normal"

$ws.Cells.Item(30, 2).Value = "Hana Murata"
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 0

# --- Row 31 : card "683506", client name not filled in yet ("a continuer") -
$ws.Cells.Item(31, 1).NumberFormat = "@"
$ws.Cells.Item(31, 1).Value = "683506"
$ws.Cells.Item(31, 1).Style = "Normal"

$ws.Cells.Item(31, 2).NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = ""
$ws.Cells.Item(31, 2).Style = "Normal"

$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 0
